# Update "想去人数" (F column) counts for two events.
# The same source rows are duplicated on sheets "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 150
    $ws.Range("F5").Value = 3059
}
